$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 edits ---

# B2: "y" -> "E", and remove the yellow highlight fill (copy format from an unstyled cell first)
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B2").Value = "E"

# E2: "hghgh" -> "sdhfsjh "
$ws.Range("E2").Value = "sdhfsjh "

# Q2: numeric value 36935 -> -32506 (keep existing date format style)
$ws.Range("Q2").Value = -32506

# AA2: 12 -> 50
$ws.Range("AA2").Value = 50

# AB2: 5000 -> 20000
$ws.Range("AB2").Value = 20000

# AV2: numeric 1 -> text "r", and apply the yellow highlight fill (copy format from M2 which has it)
$ws.Range("M2").Copy()
$ws.Range("AV2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AV2").Value = "r"

# CQ2: 1 -> 9
$ws.Range("CQ2").Value = 9

# EV2: 0 -> 2000000
$ws.Range("EV2").Value = 2000000

$excel.CutCopyMode = 0
